# chore: update Sheets via scheduled runner
#
# Refreshes the market-price-derived columns (H:N - currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) for a handful
# of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, as pulled
# by the scheduled price-sync job. Only the affected cells are touched;
# everything else in the workbook is left untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 1672.9474
$ws.Range("I86").Value = 1183.6666
$ws.Range("J86").Value = 2511.7144
$ws.Range("K86").Value = 1183.6666
$ws.Range("L86").Value = 2511.7144
$ws.Range("M86").Value = -60.66660000000002
$ws.Range("N86").Value = -4757.7144
# Row 89
$ws.Range("H89").Value = 1672.9474
$ws.Range("I89").Value = 1183.6666
$ws.Range("J89").Value = 2511.7144
$ws.Range("K89").Value = 5918.333000000001
$ws.Range("L89").Value = 12558.572
$ws.Range("M89").Value = -302.3330000000005
$ws.Range("N89").Value = -23790.572
# Row 103
$ws.Range("H103").Value = 588.44446
$ws.Range("I103").Value = 576.7273
$ws.Range("J103").Value = 596.5
$ws.Range("K103").Value = 1730.1819
$ws.Range("L103").Value = 1789.5
$ws.Range("M103").Value = -1144.1819
$ws.Range("N103").Value = -2961.5
# Row 125
$ws.Range("H125").Value = 10192995
$ws.Range("J125").Value = 11212245
$ws.Range("L125").Value = 100910205
$ws.Range("N125").Value = -100915125
# Row 138
$ws.Range("H138").Value = 9749583
$ws.Range("I138").Value = 3271443
$ws.Range("J138").Value = 12502793
$ws.Range("K138").Value = 9814329
$ws.Range("L138").Value = 37508379
$ws.Range("M138").Value = -9809189
$ws.Range("N138").Value = -37518659

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 19780.244
$ws.Range("I32").Value = 3911.5088
$ws.Range("J32").Value = 72987.17999999999
$ws.Range("K32").Value = 3911.5088
$ws.Range("L32").Value = 72987.17999999999
$ws.Range("M32").Value = -3624.5088
$ws.Range("N32").Value = -73561.17999999999
# Row 102
$ws.Range("H102").Value = 1765.8572
$ws.Range("I102").Value = 1854.6923
$ws.Range("J102").Value = 611
$ws.Range("K102").Value = 1854.6923
$ws.Range("L102").Value = 611
$ws.Range("M102").Value = -232.6922999999999
$ws.Range("N102").Value = -3855
# Row 122
$ws.Range("H122").Value = 6656.769
$ws.Range("I122").Value = 6923
$ws.Range("J122").Value = 6428.5713
$ws.Range("K122").Value = 20769
$ws.Range("L122").Value = 19285.7139
$ws.Range("M122").Value = -18319
$ws.Range("N122").Value = -24185.7139
# Row 132
$ws.Range("H132").Value = 3085.9333
$ws.Range("I132").Value = 2709.52
$ws.Range("K132").Value = 8128.559999999999
$ws.Range("M132").Value = -5598.559999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 2542.1052
$ws.Range("I99").Value = 2600
$ws.Range("J99").Value = 2515.3845
$ws.Range("K99").Value = 2600
$ws.Range("L99").Value = 2515.3845
$ws.Range("M99").Value = -1102
$ws.Range("N99").Value = -5511.3845

$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 18000
$ws.Range("J3").Value = 2000
$ws.Range("L3").Value = 2000
$ws.Range("N3").Value = -2226
# Row 31
$ws.Range("H31").Value = 5224.2
$ws.Range("I31").Value = 1778.6666
$ws.Range("J31").Value = 7891.7095
$ws.Range("K31").Value = 1778.6666
$ws.Range("L31").Value = 7891.7095
$ws.Range("M31").Value = -1483.6666
$ws.Range("N31").Value = -8481.709500000001
# Row 34
$ws.Range("H34").Value = 5224.2
$ws.Range("I34").Value = 1778.6666
$ws.Range("J34").Value = 7891.7095
$ws.Range("K34").Value = 1778.6666
$ws.Range("L34").Value = 7891.7095
$ws.Range("M34").Value = -1576.6666
$ws.Range("N34").Value = -8295.709500000001
# Row 132
$ws.Range("H132").Value = 3813.3333
$ws.Range("I132").Value = 1826.6666
$ws.Range("J132").Value = 4806.6665
$ws.Range("K132").Value = 5479.9998
$ws.Range("L132").Value = 14419.9995
$ws.Range("M132").Value = -2949.9998
$ws.Range("N132").Value = -19479.9995
# Row 134
$ws.Range("H134").Value = 2831.9048
$ws.Range("I134").Value = 1396.4117
$ws.Range("J134").Value = 8932.75
$ws.Range("K134").Value = 4189.2351
$ws.Range("L134").Value = 26798.25
$ws.Range("M134").Value = -1654.2351
$ws.Range("N134").Value = -31868.25

$ws = $wb.Worksheets.Item("CUL")
# Row 6
$ws.Range("H6").Value = 207.875
$ws.Range("I6").Value = 94.71429000000001
$ws.Range("K6").Value = 284.14287
$ws.Range("M6").Value = -171.14287
# Row 98
$ws.Range("H98").Value = 366.66666
$ws.Range("J98").Value = 366.66666
$ws.Range("L98").Value = 1099.99998
$ws.Range("N98").Value = -4095.99998
# Row 131
$ws.Range("H131").Value = 6411700
$ws.Range("J131").Value = 7408987.5
$ws.Range("L131").Value = 22226962.5
$ws.Range("N131").Value = -22237042.5
# Row 137
$ws.Range("H137").Value = 5053756.5
$ws.Range("I137").Value = 14289366
$ws.Range("J137").Value = 80735.766
$ws.Range("K137").Value = 42868098
$ws.Range("L137").Value = 242207.298
$ws.Range("M137").Value = -42862998
$ws.Range("N137").Value = -252407.298

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2968.1365
$ws.Range("I80").Value = 2943.75
$ws.Range("J80").Value = 3033.1667
$ws.Range("K80").Value = 2943.75
$ws.Range("L80").Value = 3033.1667
$ws.Range("M80").Value = -1945.75
$ws.Range("N80").Value = -5029.1667
# Row 83
$ws.Range("H83").Value = 2968.1365
$ws.Range("I83").Value = 2943.75
$ws.Range("J83").Value = 3033.1667
$ws.Range("K83").Value = 14718.75
$ws.Range("L83").Value = 15165.8335
$ws.Range("M83").Value = -9726.75
$ws.Range("N83").Value = -25149.8335
# Row 132
$ws.Range("H132").Value = 3312.2
$ws.Range("I132").Value = 3140.375
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 9421.125
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -6891.125
$ws.Range("N132").Value = -17058.5

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 1722.4117
$ws.Range("I136").Value = 1659.4546
$ws.Range("J136").Value = 3800
$ws.Range("K136").Value = 4978.3638
$ws.Range("L136").Value = 11400
$ws.Range("M136").Value = -2428.3638
$ws.Range("N136").Value = -16500

$ws = $wb.Worksheets.Item("WVR")
# Row 115
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
# Row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
# Row 118
$ws.Range("H118").Value = 27980
$ws.Range("J118").Value = 27980
$ws.Range("L118").Value = 27980
$ws.Range("N118").Value = -31294
# Row 122
$ws.Range("H122").Value = 1900
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
# Row 132
$ws.Range("H132").Value = 6043.875
$ws.Range("I132").Value = 8488.223
$ws.Range("J132").Value = 2901.1428
$ws.Range("K132").Value = 25464.669
$ws.Range("L132").Value = 8703.428400000001
$ws.Range("M132").Value = -22934.669
$ws.Range("N132").Value = -13763.4284
